$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Qwen2-5B-FoCus-length_prior"
$ws.Range("B6").Value = "0.83 ± 0.12"
$ws.Range("C6").Value = "0.48 ± 0.47"
$ws.Range("D6").Value = "-0.18 ± 0.62"
$ws.Range("E6").Value = "0.32 ± 0.68"
$ws.Range("F6").Value = "0.27 ± 0.19"
$ws.Range("G6").Value = "0.241 ± 0.00"
